$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rows 30-32: larval number (column F) corrected to 12 for all three Control
# replicates.
# ---------------------------------------------------------------------------
$ws.Range("F30").Value = 12
$ws.Range("F31").Value = 12
$ws.Range("F32").Value = 12

# ---------------------------------------------------------------------------
# Row 33: becomes the pooled "Control" tube (chambers 1,2,3 combined) instead
# of the first "Medium" replicate. The Squarical cell (D33) is removed and a
# new Chamber cell (E33) holding "1,2,3" is introduced.
# ---------------------------------------------------------------------------
$ws.Range("C33").Value = "Control"
$ws.Range("D33").Clear()
$ws.Range("E33").Value = "1,2,3"
$ws.Range("F33").Value = 10
$ws.Range("G33").Value = "PC2a"
# G33 used to carry the "Tube ID" style (12); it now matches the plain
# centered style used elsewhere in the row (copy format from F33, style 11).
$ws.Range("F33").Copy() | Out-Null
$ws.Range("G33").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# Rows 34-35: the "Medium" replicates shift up one Squarical/Larval slot.
# ---------------------------------------------------------------------------
$ws.Range("D34").Value = 1
$ws.Range("F34").Value = 12
$ws.Range("G34").Value = "PM1s"

$ws.Range("D35").Value = 2
$ws.Range("F35").Value = 12
$ws.Range("G35").Value = "PM2s"

# ---------------------------------------------------------------------------
# Row 36: becomes the third "Medium" replicate (was the first "High" one).
# ---------------------------------------------------------------------------
$ws.Range("C36").Value = "Medium"
$ws.Range("D36").Value = 3
$ws.Range("G36").Value = "PM3s"
$ws.Range("L36").ClearContents()

# ---------------------------------------------------------------------------
# Row 37: becomes the pooled "Medium" tube (chambers 1,2,3 combined).
# ---------------------------------------------------------------------------
$ws.Range("C37").Value = "Medium"
$ws.Range("D37").Clear()
$ws.Range("E37").Value = "1,2,3"
$ws.Range("F37").Value = 9
$ws.Range("G37").Value = "PM1a"
$ws.Range("F37").Copy() | Out-Null
$ws.Range("G37").PasteSpecial(-4122) | Out-Null
$ws.Range("L37").ClearContents()

# ---------------------------------------------------------------------------
# Row 38: first "High" replicate (was the third one).
# ---------------------------------------------------------------------------
$ws.Range("D38").Value = 1
$ws.Range("F38").Value = 12
$ws.Range("G38").Value = "PH2a"
$ws.Range("L38").ClearContents()

# ---------------------------------------------------------------------------
# Rows 39-40: newly filled-in "High" replicates 2 and 3 (previously just a
# date stamp in column A). Build them out to match the rest of the table,
# copying formatting from the now-complete row 38.
# ---------------------------------------------------------------------------
$ws.Range("C39").Value = "High"
$ws.Range("D39").Value = 2
$ws.Range("F39").Value = 12
$ws.Range("G39").Value = "PH3a"
$ws.Range("H39").Value = "RNA-DNA shield"
$ws.Range("I39").Value = 300
$ws.Range("J39").Value = "-80C"

$ws.Range("C40").Value = "High"
$ws.Range("D40").Value = 3
$ws.Range("F40").Value = 11
$ws.Range("G40").Value = "PH1b"
$ws.Range("H40").Value = "RNA-DNA shield"
$ws.Range("I40").Value = 300
$ws.Range("J40").Value = "-80C"

# Copy the row-38 formatting across so the new cells pick up the correct
# styles (centered numbers, Tube-ID styling, Notes styling, etc.).
$ws.Range("A38:L38").Copy() | Out-Null
$ws.Range("A39:L39").PasteSpecial(-4122) | Out-Null
$ws.Range("A38:L38").Copy() | Out-Null
$ws.Range("A40:L40").PasteSpecial(-4122) | Out-Null

# The two new rows have no Notes entry, so make sure L39/L40 are blank.
$ws.Range("L39").ClearContents()
$ws.Range("L40").ClearContents()

# ---------------------------------------------------------------------------
# Rows 41-44 no longer hold the placeholder date stamp.
# ---------------------------------------------------------------------------
$ws.Range("A41").Clear()
$ws.Range("A42").Clear()
$ws.Range("A43").Clear()
$ws.Range("A44").Clear()

# ---------------------------------------------------------------------------
# Leave the selection where the author's last edit was.
# ---------------------------------------------------------------------------
$ws.Range("E37").Select()
